$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cellerrors")

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Formula = "=10+10"
$ws.Range("C8").Formula = "=SUM(C2:C7)"

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Hello"
$ws.Range("C9").Formula = "=LOWER(B9)"

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "abcd"
$ws.Range("C10").Value = "wxyz"

$ws.Activate()
$ws.Range("C11").Select()
